# The four observation records currently on rows 17-20 get shuffled: each
# record's row-varying fields move up one row (18->17, 19->18, 20->19), and
# the record that was on row 17 wraps around to row 20.  Columns that are
# identical across all four rows (P, S, T, U, V, W, Y, Z, AA, AB, AD, AE,
# AG, AT, AW, AX, AY) are left untouched.
#
# Target (post-edit) values per row, taken from the diff:
#
# Row 17: Id=130960789 sortOrd=79243 TaxonId=6425   Artnamn=Garnlav            VetNamn=Alectoria sarmentosa  Auktor=(Ach.) Ach.                      Aktivitet=(none)      Ost=446284 Nord=6759886 PublikKommentar=(none)
# Row 18: Id=130960843 sortOrd=79243 TaxonId=6425   Artnamn=Garnlav            VetNamn=Alectoria sarmentosa  Auktor=(Ach.) Ach.                      Aktivitet=(none)      Ost=446247 Nord=6759903 PublikKommentar=(none)
# Row 19: Id=130961956 sortOrd=79862 TaxonId=6453   Artnamn=Vedskivlav         VetNamn=Hertelidea botryosa   Auktor=(Fr.) Printzen & Kantvilas        Aktivitet=(none)      Ost=446084 Nord=6759981 PublikKommentar=Miljöbilder
# Row 20: Id=130960378 sortOrd=57884 TaxonId=100109 Artnamn=Tretåig hackspett  VetNamn=Picoides tridactylus  Auktor=(Linnaeus, 1758)                  Aktivitet=äldre spår  Ost=446272 Nord=6759739 PublikKommentar=(none)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17 (was row 18's data) ---
$ws.Range("A17").Value = 130960789
$ws.Range("B17").Value = 79243
$ws.Range("E17").Value = 6425
$ws.Range("F17").Value = "Garnlav"
$ws.Range("G17").Value = "Alectoria sarmentosa"
$ws.Range("H17").Value = "(Ach.) Ach."
$ws.Range("M17").ClearContents()
$ws.Range("Q17").Value = 446284
$ws.Range("R17").Value = 6759886
$ws.Range("AC17").ClearContents()

# --- Row 18 (was row 19's data) ---
$ws.Range("A18").Value = 130960843
$ws.Range("B18").Value = 79243
$ws.Range("E18").Value = 6425
$ws.Range("F18").Value = "Garnlav"
$ws.Range("G18").Value = "Alectoria sarmentosa"
$ws.Range("H18").Value = "(Ach.) Ach."
$ws.Range("M18").ClearContents()
$ws.Range("Q18").Value = 446247
$ws.Range("R18").Value = 6759903
$ws.Range("AC18").ClearContents()

# --- Row 19 (was row 20's data) ---
$ws.Range("A19").Value = 130961956
$ws.Range("B19").Value = 79862
$ws.Range("E19").Value = 6453
$ws.Range("F19").Value = "Vedskivlav"
$ws.Range("G19").Value = "Hertelidea botryosa"
$ws.Range("H19").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("M19").ClearContents()
$ws.Range("Q19").Value = 446084
$ws.Range("R19").Value = 6759981
$ws.Range("AC19").Value = "Miljöbilder"

# --- Row 20 (was row 17's data) ---
$ws.Range("A20").Value = 130960378
$ws.Range("B20").Value = 57884
$ws.Range("E20").Value = 100109
$ws.Range("F20").Value = "Tretåig hackspett"
$ws.Range("G20").Value = "Picoides tridactylus"
$ws.Range("H20").Value = "(Linnaeus, 1758)"
$ws.Range("M20").Value = "äldre spår"
$ws.Range("Q20").Value = 446272
$ws.Range("R20").Value = 6759739
$ws.Range("AC20").ClearContents()
